$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "name, email address" placeholder to the client's email.
$ws.Range("A31").Value = "charlie.charlie@mail.com"

# Apply the client discount amount (Client discount row, Unit Price column).
$ws.Range("E18").Value = 100
